$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 96.62005599999999
$ws.Range("H2").Value = 289.860168
$ws.Range("I2").Value = 0.2116037895476247
$ws.Range("J2").Value = 0.2183905833651517
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 11.465689
$ws.Range("N2").Value = 34.397067
$ws.Range("O2").Value = 0.1125836808441207
$ws.Range("P2").Value = 0.1279391038575984
$ws.Range("Q2").Value = 1107.815513258584
$ws.Range("R2").Value = 9970.339619327255
$ws.Range("S2").Value = 0.02382313350783626
$ws.Range("T2").Value = 0.02794069552667565

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 96.62005599999999
$ws.Range("H3").Value = 289.860168
$ws.Range("I3").Value = 0.2116037895476247
$ws.Range("J3").Value = 0.2183905833651517
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.399706666666667
$ws.Range("N3").Value = 4.199120000000001
$ws.Range("O3").Value = 0.01374397374945266
$ws.Range("P3").Value = 0.01561853078317749
$ws.Range("Q3").Value = 135.2397365169067
$ws.Range("R3").Value = 1217.15762865216
$ws.Range("S3").Value = 0.002908276928827258
$ws.Range("T3").Value = 0.003410940049044712

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 96.62005599999999
$ws.Range("H4").Value = 289.860168
$ws.Range("I4").Value = 0.2116037895476247
$ws.Range("J4").Value = 0.2183905833651517
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 26.364677
$ws.Range("N4").Value = 79.094031
$ws.Range("O4").Value = 0.2588795475724425
$ws.Range("P4").Value = 0.2941884389917637
$ws.Range("Q4").Value = 2547.356568161912
$ws.Range("R4").Value = 22926.20911345721
$ws.Range("S4").Value = 0.0547798933027034
$ws.Range("T4").Value = 0.06424798481069462

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 96.62005599999999
$ws.Range("H5").Value = 289.860168
$ws.Range("I5").Value = 0.2116037895476247
$ws.Range("J5").Value = 0.2183905833651517
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.94195566666667
$ws.Range("N5").Value = 77.825867
$ws.Range("O5").Value = 0.2547287700938277
$ws.Range("P5").Value = 0.2894715320036049
$ws.Range("Q5").Value = 2506.513209262851
$ws.Range("R5").Value = 22558.61888336566
$ws.Range("S5").Value = 0.05390157305865959
$ws.Range("T5").Value = 0.06321785674187146

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 96.62005599999999
$ws.Range("H6").Value = 289.860168
$ws.Range("I6").Value = 0.2116037895476247
$ws.Range("J6").Value = 0.2183905833651517
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 36.6694545
$ws.Range("N6").Value = 73.338909
$ws.Range("O6").Value = 0.3600640277401564
$ws.Range("P6").Value = 0.2727823943638554
$ws.Range("Q6").Value = 3543.004747279452
$ws.Range("R6").Value = 21258.02848367671
$ws.Range("S6").Value = 0.07619091274959815
$ws.Range("T6").Value = 0.05957310623686526

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 122.3539896666667
$ws.Range("H7").Value = 367.061969
$ws.Range("I7").Value = 0.2679626668787852
$ws.Range("J7").Value = 0.2765570657541026
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 11.465689
$ws.Range("N7").Value = 34.397067
$ws.Range("O7").Value = 0.1125836808441207
$ws.Range("P7").Value = 0.1279391038575984
$ws.Range("Q7").Value = 1402.872793427214
$ws.Range("R7").Value = 12625.85514084492
$ws.Range("S7").Value = 0.03016822336602059
$ws.Range("T7").Value = 0.03538246315806681

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 122.3539896666667
$ws.Range("H8").Value = 367.061969
$ws.Range("I8").Value = 0.2679626668787852
$ws.Range("J8").Value = 0.2765570657541026
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.399706666666667
$ws.Range("N8").Value = 4.199120000000001
$ws.Range("O8").Value = 0.01374397374945266
$ws.Range("P8").Value = 0.01561853078317749
$ws.Range("Q8").Value = 171.2596950296978
$ws.Range("R8").Value = 1541.33725526728
$ws.Range("S8").Value = 0.003682871859415352
$ws.Range("T8").Value = 0.004319415044785693

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 122.3539896666667
$ws.Range("H9").Value = 367.061969
$ws.Range("I9").Value = 0.2679626668787852
$ws.Range("J9").Value = 0.2765570657541026
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 26.364677
$ws.Range("N9").Value = 79.094031
$ws.Range("O9").Value = 0.2588795475724425
$ws.Range("P9").Value = 0.2941884389917637
$ws.Range("Q9").Value = 3225.823417223004
$ws.Range("R9").Value = 29032.41075500704
$ws.Range("S9").Value = 0.06937005396788504
$ws.Range("T9").Value = 0.08135989146634198

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 122.3539896666667
$ws.Range("H10").Value = 367.061969
$ws.Range("I10").Value = 0.2679626668787852
$ws.Range("J10").Value = 0.2765570657541026
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 25.94195566666667
$ws.Range("N10").Value = 77.825867
$ws.Range("O10").Value = 0.2547287700938277
$ws.Range("P10").Value = 0.2894715320036049
$ws.Range("Q10").Value = 3174.101775572458
$ws.Range("R10").Value = 28566.91598015212
$ws.Range("S10").Value = 0.06825780056509502
$ws.Range("T10").Value = 0.08005539751026176

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 122.3539896666667
$ws.Range("H11").Value = 367.061969
$ws.Range("I11").Value = 0.2679626668787852
$ws.Range("J11").Value = 0.2765570657541026
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 36.6694545
$ws.Range("N11").Value = 73.338909
$ws.Range("O11").Value = 0.3600640277401564
$ws.Range("P11").Value = 0.2727823943638554
$ws.Range("Q11").Value = 4486.654056975303
$ws.Range("R11").Value = 26919.92434185182
$ws.Range("S11").Value = 0.09648371712036923
$ws.Range("T11").Value = 0.07543989857464631

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 83.74384566666667
$ws.Range("H12").Value = 251.231537
$ws.Range("I12").Value = 0.1834041070557659
$ws.Range("J12").Value = 0.1892864490617203
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 11.465689
$ws.Range("N12").Value = 34.397067
$ws.Range("O12").Value = 0.1125836808441207
$ws.Range("P12").Value = 0.1279391038575984
$ws.Range("Q12").Value = 960.1808900779977
$ws.Range("R12").Value = 8641.628010701979
$ws.Range("S12").Value = 0.02064830945426729
$ws.Range("T12").Value = 0.02421713866534345

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 83.74384566666667
$ws.Range("H13").Value = 251.231537
$ws.Range("I13").Value = 0.1834041070557659
$ws.Range("J13").Value = 0.1892864490617203
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.399706666666667
$ws.Range("N13").Value = 4.199120000000001
$ws.Range("O13").Value = 0.01374397374945266
$ws.Range("P13").Value = 0.01561853078317749
$ws.Range("Q13").Value = 117.2168190719378
$ws.Range("R13").Value = 1054.95137164744
$ws.Range("S13").Value = 0.002520701232916252
$ws.Range("T13").Value = 0.002956376231508837

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 83.74384566666667
$ws.Range("H14").Value = 251.231537
$ws.Range("I14").Value = 0.1834041070557659
$ws.Range("J14").Value = 0.1892864490617203
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 26.364677
$ws.Range("N14").Value = 79.094031
$ws.Range("O14").Value = 0.2588795475724425
$ws.Range("P14").Value = 0.2941884389917637
$ws.Range("Q14").Value = 2207.879441739517
$ws.Range("R14").Value = 19870.91497565565
$ws.Range("S14").Value = 0.04747957225752447
$ws.Range("T14").Value = 0.0556858849717615

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 83.74384566666667
$ws.Range("H15").Value = 251.231537
$ws.Range("I15").Value = 0.1834041070557659
$ws.Range("J15").Value = 0.1892864490617203
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 25.94195566666667
$ws.Range("N15").Value = 77.825867
$ws.Range("O15").Value = 0.2547287700938277
$ws.Range("P15").Value = 0.2894715320036049
$ws.Range("Q15").Value = 2172.479131640842
$ws.Range("R15").Value = 19552.31218476758
$ws.Range("S15").Value = 0.04671830262047195
$ws.Range("T15").Value = 0.0547930383974185

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 83.74384566666667
$ws.Range("H16").Value = 251.231537
$ws.Range("I16").Value = 0.1834041070557659
$ws.Range("J16").Value = 0.1892864490617203
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 36.6694545
$ws.Range("N16").Value = 73.338909
$ws.Range("O16").Value = 0.3600640277401564
$ws.Range("P16").Value = 0.2727823943638554
$ws.Range("Q16").Value = 3070.841138328856
$ws.Range("R16").Value = 18425.04682997313
$ws.Range("S16").Value = 0.06603722149058591
$ws.Range("T16").Value = 0.05163401079568803

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 111.321218
$ws.Range("H17").Value = 333.963654
$ws.Range("I17").Value = 0.2438002270031519
$ws.Range("J17").Value = 0.2516196610353779
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 11.465689
$ws.Range("N17").Value = 34.397067
$ws.Range("O17").Value = 0.1125836808441207
$ws.Range("P17").Value = 0.1279391038575984
$ws.Range("Q17").Value = 1276.374464689202
$ws.Range("R17").Value = 11487.37018220282
$ws.Range("S17").Value = 0.02744792694664703
$ws.Range("T17").Value = 0.03219199394581892

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 111.321218
$ws.Range("H18").Value = 333.963654
$ws.Range("I18").Value = 0.2438002270031519
$ws.Range("J18").Value = 0.2516196610353779
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 1.399706666666667
$ws.Range("N18").Value = 4.199120000000001
$ws.Range("O18").Value = 0.01374397374945266
$ws.Range("P18").Value = 0.01561853078317749
$ws.Range("Q18").Value = 155.8170509760534
$ws.Range("R18").Value = 1402.35345878448
$ws.Range("S18").Value = 0.003350783920041919
$ws.Range("T18").Value = 0.003929929421533735

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 111.321218
$ws.Range("H19").Value = 333.963654
$ws.Range("I19").Value = 0.2438002270031519
$ws.Range("J19").Value = 0.2516196610353779
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 26.364677
$ws.Range("N19").Value = 79.094031
$ws.Range("O19").Value = 0.2588795475724425
$ws.Range("P19").Value = 0.2941884389917637
$ws.Range("Q19").Value = 2934.947955816586
$ws.Range("R19").Value = 26414.53160234927
$ws.Range("S19").Value = 0.06311489246463473
$ws.Range("T19").Value = 0.07402359529963451

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 111.321218
$ws.Range("H20").Value = 333.963654
$ws.Range("I20").Value = 0.2438002270031519
$ws.Range("J20").Value = 0.2516196610353779
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 25.94195566666667
$ws.Range("N20").Value = 77.825867
$ws.Range("O20").Value = 0.2547287700938277
$ws.Range("P20").Value = 0.2894715320036049
$ws.Range("Q20").Value = 2887.890102115336
$ws.Range("R20").Value = 25991.01091903802
$ws.Range("S20").Value = 0.06210293197310889
$ws.Range("T20").Value = 0.07283672876213859

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 111.321218
$ws.Range("H21").Value = 333.963654
$ws.Range("I21").Value = 0.2438002270031519
$ws.Range("J21").Value = 0.2516196610353779
$ws.Range("K21").Value = 2
$ws.Range("M21").Value = 36.6694545
$ws.Range("N21").Value = 73.338909
$ws.Range("O21").Value = 0.3600640277401564
$ws.Range("P21").Value = 0.2727823943638554
$ws.Range("Q21").Value = 4082.088338335581
$ws.Range("R21").Value = 24492.53003001349
$ws.Range("S21").Value = 0.08778369169871932
$ws.Range("T21").Value = 0.06863741360625207

$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 42.569235
$ws.Range("H22").Value = 85.13847
$ws.Range("I22").Value = 0.09322920951467238
$ws.Range("J22").Value = 0.06414624078364733
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 11.465689
$ws.Range("N22").Value = 34.397067
$ws.Range("O22").Value = 0.1125836808441207
$ws.Range("P22").Value = 0.1279391038575984
$ws.Range("Q22").Value = 488.0856094779149
$ws.Range("R22").Value = 2928.51365686749
$ws.Range("S22").Value = 0.01049608756934954
$ws.Range("T22").Value = 0.008206812561693571

$ws.Range("E23").Value = 2
$ws.Range("G23").Value = 42.569235
$ws.Range("H23").Value = 85.13847
$ws.Range("I23").Value = 0.09322920951467238
$ws.Range("J23").Value = 0.06414624078364733
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 1.399706666666667
$ws.Range("N23").Value = 4.199120000000001
$ws.Range("O23").Value = 0.01374397374945266
$ws.Range("P23").Value = 0.01561853078317749
$ws.Range("Q23").Value = 59.5844420244
$ws.Range("R23").Value = 357.5066521464
$ws.Range("S23").Value = 0.001281339808251879
$ws.Range("T23").Value = 0.001001870036304511

$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 42.569235
$ws.Range("H24").Value = 85.13847
$ws.Range("I24").Value = 0.09322920951467238
$ws.Range("J24").Value = 0.06414624078364733
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 26.364677
$ws.Range("N24").Value = 79.094031
$ws.Range("O24").Value = 0.2588795475724425
$ws.Range("P24").Value = 0.2941884389917637
$ws.Range("Q24").Value = 1122.324130912095
$ws.Range("R24").Value = 6733.944785472569
$ws.Range("S24").Value = 0.02413513557969484
$ws.Range("T24").Value = 0.01887108244333102

$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 42.569235
$ws.Range("H25").Value = 85.13847
$ws.Range("I25").Value = 0.09322920951467238
$ws.Range("J25").Value = 0.06414624078364733
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 25.94195566666667
$ws.Range("N25").Value = 77.825867
$ws.Range("O25").Value = 0.2547287700938277
$ws.Range("P25").Value = 0.2894715320036049
$ws.Range("Q25").Value = 1104.329207133915
$ws.Range("R25").Value = 6625.97524280349
$ws.Range("S25").Value = 0.02374816187649227
$ws.Range("T25").Value = 0.01856851059191451

$ws.Range("E26").Value = 2
$ws.Range("G26").Value = 42.569235
$ws.Range("H26").Value = 85.13847
$ws.Range("I26").Value = 0.09322920951467238
$ws.Range("J26").Value = 0.06414624078364733
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 36.6694545
$ws.Range("N26").Value = 73.338909
$ws.Range("O26").Value = 0.3600640277401564
$ws.Range("P26").Value = 0.2727823943638554
$ws.Range("Q26").Value = 1560.990625932307
$ws.Range("R26").Value = 6243.96250372923
$ws.Range("S26").Value = 0.03356848468088385
$ws.Range("T26").Value = 0.01749796515040371
